$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2140.6
$ws.Range("I2").Value = 289.8
$ws.Range("K2").Value = 289.8
$ws.Range("M2").Value = -176.8
$ws.Range("H43").Value = 3740.2
$ws.Range("I43").Value = 2000
$ws.Range("J43").Value = 4175.25
$ws.Range("K43").Value = 2000
$ws.Range("L43").Value = 4175.25
$ws.Range("M43").Value = -1931
$ws.Range("N43").Value = -4313.25
$ws.Range("H62").Value = 8185.1143
$ws.Range("I62").Value = 8185.1143
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 8185.1143
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -7561.1143
$ws.Range("N62").Value = $null
$ws.Range("H65").Value = 8185.1143
$ws.Range("I65").Value = 8185.1143
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 40925.5715
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -37805.5715
$ws.Range("N65").Value = $null
$ws.Range("H132").Value = 3132
$ws.Range("I132").Value = 3132
$ws.Range("K132").Value = 9396
$ws.Range("M132").Value = -6866
$ws.Range("H137").Value = 4530.727
$ws.Range("I137").Value = 4027.423
$ws.Range("K137").Value = 12082.269
$ws.Range("M137").Value = -9532.269
$ws.Range("H138").Value = 7316.1113
$ws.Range("J138").Value = 7223.5
$ws.Range("L138").Value = 21670.5
$ws.Range("N138").Value = -31950.5
# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2652.9333
$ws.Range("I2").Value = 2263.389
$ws.Range("J2").Value = 4211.1113
$ws.Range("K2").Value = 2263.389
$ws.Range("L2").Value = 4211.1113
$ws.Range("M2").Value = -2150.389
$ws.Range("N2").Value = -4437.1113
$ws.Range("H44").Value = 69158
$ws.Range("J44").Value = 69158
$ws.Range("L44").Value = 69158
$ws.Range("N44").Value = -70134
$ws.Range("H55").Value = 49919.4
$ws.Range("J55").Value = 49919.4
$ws.Range("L55").Value = 49919.4
$ws.Range("N55").Value = -50549.4
$ws.Range("H61").Value = 20014
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").Value = $null
$ws.Range("H74").Value = 7234.6924
$ws.Range("I74").Value = 6143.375
$ws.Range("K74").Value = 6143.375
$ws.Range("M74").Value = -5269.375
$ws.Range("H77").Value = 7234.6924
$ws.Range("I77").Value = 6143.375
$ws.Range("K77").Value = 30716.875
$ws.Range("M77").Value = -26348.875
$ws.Range("H88").Value = 1970.7142
$ws.Range("I88").Value = 1665
$ws.Range("K88").Value = 1665
$ws.Range("M88").Value = -1259
$ws.Range("H91").Value = 1970.7142
$ws.Range("I91").Value = 1665
$ws.Range("K91").Value = 1665
$ws.Range("M91").Value = -261
$ws.Range("H116").Value = 2652.9333
$ws.Range("I116").Value = 2263.389
$ws.Range("J116").Value = 4211.1113
$ws.Range("K116").Value = 2263.389
$ws.Range("L116").Value = 4211.1113
$ws.Range("M116").Value = 30.61099999999988
$ws.Range("N116").Value = -8799.1113
$ws.Range("H132").Value = 3478.6667
$ws.Range("I132").Value = 2996.6223
$ws.Range("J132").Value = 5888.8887
$ws.Range("K132").Value = 8989.866900000001
$ws.Range("L132").Value = 17666.6661
$ws.Range("M132").Value = -6459.866900000001
$ws.Range("N132").Value = -22726.6661
$ws.Range("H136").Value = 20014
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").Value = $null
# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2652.9333
$ws.Range("I3").Value = 2263.389
$ws.Range("J3").Value = 4211.1113
$ws.Range("K3").Value = 2263.389
$ws.Range("L3").Value = 4211.1113
$ws.Range("M3").Value = -2149.389
$ws.Range("N3").Value = -4439.1113
$ws.Range("H86").Value = 16866.223
$ws.Range("J86").Value = 12099.2
$ws.Range("L86").Value = 12099.2
$ws.Range("N86").Value = -14345.2
$ws.Range("H89").Value = 16866.223
$ws.Range("J89").Value = 12099.2
$ws.Range("L89").Value = 60496
$ws.Range("N89").Value = -71728
$ws.Range("H134").Value = 6489.6387
$ws.Range("I134").Value = 5567.963
$ws.Range("K134").Value = 16703.889
$ws.Range("M134").Value = -14168.889
# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 89001
$ws.Range("I23").Value = 89001
$ws.Range("K23").Value = 89001
$ws.Range("M23").Value = -88761
$ws.Range("H27").Value = 89001
$ws.Range("I27").Value = 89001
$ws.Range("K27").Value = 89001
$ws.Range("M27").Value = -88809
$ws.Range("H31").Value = 16898
$ws.Range("I31").Value = 7082.857
$ws.Range("K31").Value = 7082.857
$ws.Range("M31").Value = -6787.857
$ws.Range("H34").Value = 16898
$ws.Range("I34").Value = 7082.857
$ws.Range("K34").Value = 7082.857
$ws.Range("M34").Value = -6880.857
$ws.Range("H36").Value = 13833.333
$ws.Range("I36").Value = 8500
$ws.Range("J36").Value = 16500
$ws.Range("K36").Value = 8500
$ws.Range("L36").Value = 16500
$ws.Range("M36").Value = -8112
$ws.Range("N36").Value = -17276
$ws.Range("H40").Value = 13833.333
$ws.Range("I40").Value = 8500
$ws.Range("J40").Value = 16500
$ws.Range("K40").Value = 8500
$ws.Range("L40").Value = 16500
$ws.Range("M40").Value = -8340
$ws.Range("N40").Value = -16820
$ws.Range("H54").Value = 39624.125
$ws.Range("J54").Value = 39624.125
$ws.Range("L54").Value = 39624.125
$ws.Range("N54").Value = -40940.125
$ws.Range("H55").Value = 19344.666
$ws.Range("I55").Value = 25000
$ws.Range("J55").Value = 16517
$ws.Range("K55").Value = 25000
$ws.Range("L55").Value = 16517
$ws.Range("M55").Value = -24685
$ws.Range("N55").Value = -17147
$ws.Range("H56").Value = 20000
$ws.Range("J56").Value = 20000
$ws.Range("L56").Value = 20000
$ws.Range("N56").Value = -21690
$ws.Range("H57").Value = 55499.25
$ws.Range("J57").Value = 57332.332
$ws.Range("L57").Value = 57332.332
$ws.Range("N57").Value = -58452.332
$ws.Range("H64").Value = 49998.5
$ws.Range("J64").Value = 49998.5
$ws.Range("L64").Value = 49998.5
$ws.Range("N64").Value = -50494.5
$ws.Range("H67").Value = 49998.5
$ws.Range("J67").Value = 49998.5
$ws.Range("L67").Value = 49998.5
$ws.Range("N67").Value = -51714.5
$ws.Range("H132").Value = 5052.469
$ws.Range("I132").Value = 4604
$ws.Range("K132").Value = 13812
$ws.Range("M132").Value = -11282
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = $null
# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 190
$ws.Range("I2").Value = 280
$ws.Range("J2").Value = 40
$ws.Range("K2").Value = 1680
$ws.Range("L2").Value = 240
$ws.Range("M2").Value = -1567
$ws.Range("N2").Value = -466
$ws.Range("H38").Value = 13.166667
$ws.Range("J38").Value = 11.666667
$ws.Range("L38").Value = 35.000001
$ws.Range("N38").Value = -729.000001
$ws.Range("H39").Value = 19999
$ws.Range("J39").Value = 19999
$ws.Range("L39").Value = 59997
$ws.Range("N39").Value = -60585
$ws.Range("H113").Value = 1795.8334
$ws.Range("I113").Value = 1795.25
$ws.Range("K113").Value = 5385.75
$ws.Range("M113").Value = -3215.75
$ws.Range("H129").Value = 8774949
$ws.Range("I129").Value = 1237.4286
$ws.Range("J129").Value = 33341340
$ws.Range("K129").Value = 3712.2858
$ws.Range("L129").Value = 100024020
$ws.Range("M129").Value = 1287.7142
$ws.Range("N129").Value = -100034020
$ws.Range("H133").Value = 9201.6
$ws.Range("I133").Value = 6258
$ws.Range("J133").Value = 9937.5
$ws.Range("K133").Value = 18774
$ws.Range("L133").Value = 29812.5
$ws.Range("M133").Value = -13714
$ws.Range("N133").Value = -39932.5
$ws.Range("H136").Value = 3174.875
$ws.Range("I136").Value = 3174.875
$ws.Range("K136").Value = 9524.625
$ws.Range("M136").Value = -4424.625
# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 9194.409
$ws.Range("I132").Value = 8626.666999999999
$ws.Range("J132").Value = 11749.25
$ws.Range("K132").Value = 25880.001
$ws.Range("L132").Value = 35247.75
$ws.Range("M132").Value = -23350.001
$ws.Range("N132").Value = -40307.75
# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4361.619
$ws.Range("I40").Value = 3921.889
$ws.Range("K40").Value = 3921.889
$ws.Range("M40").Value = -3785.889
$ws.Range("H132").Value = 4488.125
$ws.Range("I132").Value = 3786.8518
$ws.Range("J132").Value = 8275
$ws.Range("K132").Value = 11360.5554
$ws.Range("L132").Value = 24825
$ws.Range("M132").Value = -8830.555399999999
$ws.Range("N132").Value = -29885
# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 423.75
$ws.Range("I55").Value = 500
$ws.Range("J55").Value = 398.33334
$ws.Range("K55").Value = 500
$ws.Range("L55").Value = 398.33334
$ws.Range("M55").Value = -223
$ws.Range("N55").Value = -952.33334
$ws.Range("H81").Value = 2457.6316
$ws.Range("J81").Value = 4640
$ws.Range("L81").Value = 9280
$ws.Range("N81").Value = -11402
$ws.Range("H84").Value = 2457.6316
$ws.Range("J84").Value = 4640
$ws.Range("L84").Value = 46400
$ws.Range("N84").Value = -57008
$ws.Range("H126").Value = 1940.8959
$ws.Range("I126").Value = 1547.8292
$ws.Range("K126").Value = 4643.487599999999
$ws.Range("M126").Value = -2173.487599999999
$ws.Range("H132").Value = 6728.2896
$ws.Range("I132").Value = 6493.4116
$ws.Range("K132").Value = 19480.2348
$ws.Range("M132").Value = -16950.2348
